$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metrics")

# Update the "covidlive.com.au" values (column D) for rows 63-68
$ws.Range("D63").Value = 620
$ws.Range("D64").Value = 630
$ws.Range("D65").Value = 640
$ws.Range("D66").Value = 650
$ws.Range("D67").Value = 660
$ws.Range("D68").Value = 670

# Update the selected cell/range on the sheet to D60
$ws.Range("D60").Select()
